$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "№ Задача - Заявки, включващи обобщаващи функции (SUM, COUNT, AVG, MIN, MAX)"
# column is being removed. Shift the two headers that followed it one column to
# the left (D and E) and clear out the now-unused last header cell (F).
$ws.Range("D1").Value = "№ Задача - Заявки, включващи GROUP BY, HAVING"
$ws.Range("E1").Value = "№ Задача - Заявки, включващи външни съединения и/или подзаявки"
$ws.Range("F1").Clear()

# Populate the previously empty task-number columns (C, D, E) for every student row.
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 2

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 6

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 10

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 11

$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 9

$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 13

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 5

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 7

$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = 4

$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 3

$ws.Range("C13").Value = 9
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 8

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 15

# Shrink the table to match the new A1:E14 layout (drops the removed column).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E14"))

# The last remaining column got noticeably wider after the column removal.
$ws.Columns("E").ColumnWidth = 68.45

# The active selection moved up now that there is one column fewer.
$ws.Range("A16").Select()
